# "all important files in one folder"
# - Add a qualifier note "(?)" next to the entries in rows 17 and 22 (column C)
# - Add a qualifier note "(Pierson syndrome)" next to the entry in row 22 (column E)
# - Re-apply the "duplicate values" conditional formatting on A1:A25 so a
#   fresh (identical) highlight-fill format is registered and the rule now
#   points at it
# - Move the active selection to E17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = "(?)"
$ws.Range("C22").Value = "(?)"
$ws.Range("E22").Value = "(Pierson syndrome)"

$range = $ws.Range("A1:A25")
$dupeColor = $range.FormatConditions.Item(1).Interior.Color
$range.FormatConditions.Delete()
$fc = $range.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Interior.Color = $dupeColor

[void]$ws.Range("E17").Select()
